$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zhenyun Lyu")

# Row 6 (Sprint1 "Done" - first week)
$ws.Range("A6").Value = "followed the instructions on blackboard to clone, run and test the code"
$ws.Range("B6").Value = "get more familiar with the souce code"
$ws.Range("C6").Value = "get full understand of the tasks in sprint1"
$ws.Range("D6").Value = "finished the test for the calculate_first_last function"
$ws.Range("E6").Value = "Sprint Review / retrospective"

# Row 7
$ws.Range("A7").Value = "went through the source code"
$ws.Range("B7").Clear()
$ws.Range("C7").Value = "analysed the requirements"
$ws.Range("E7").Value = "Commited task"

# Row 8
$ws.Range("A8").Value = "discussed the sprint1 with team menbers and got the divided work"

# Row 9
$ws.Range("A9").Value = "fixed the existed errors and passed the tests"

# Row 14 (Sprint1 "To do" - second week)
$ws.Range("A14").Value = "get more familiar with the souce code"
$ws.Range("B14").Value = "get full understand of the tasks in sprint1"
$ws.Range("C14").Value = "finished the test for the calculate_first_last function"

# Row 15
$ws.Range("A15").Value = "get full understand of the tasks in sprint1"
$ws.Range("B15").Value = "analysed the requirements"

# Row 16
$ws.Range("A16").Value = "analysed the requirements"
$ws.Range("B16").Value = "finish the test for the new function: calculate_first_last"

# Row 17
$ws.Range("A17").Value = "finish the test for the new function: calculate_first_last"

# Column widths for this sheet (closest achievable values given the
# runtime's character-width/pixel quantization)
$ws.Columns.Item(1).ColumnWidth = 56.714285714285715
$ws.Columns.Item(2).ColumnWidth = 43.142857142857146
$ws.Columns.Item(3).ColumnWidth = 40.57142857142857
$ws.Columns.Item(4).ColumnWidth = 41.57142857142857
$ws.Columns.Item(5).ColumnWidth = 22.714285714285715

# Make this sheet the active/selected one, with A4:G4 selected (active cell A4)
$ws.Select()
$ws.Range("A4:G4").Select()
